$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the value in A1 (test1 -> test11); A2/A5 keep their existing text
$ws.Range("A1").Value = "test11"

# Move the active selection to A2 and reset the scrolled top-left cell to A1
$ws.Application.Goto($ws.Range("A1"), $true)
$ws.Range("A2").Select()
